$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the last four columns (J,K,L,M) of the header block:
# Before: J=Comment, K=SamplePortion, L=SamplePortionUnit, M=ResultUnit
# After:  J=SamplePortion, K=SamplePortionUnit, L=ResultUnit, M=Comment
# (the "Comment" column moves to the end of the table)

$ws.Range("J1").Value = "SamplePortion"
$ws.Range("K1").Value = "SamplePortionUnit"
$ws.Range("L1").Value = "ResultUnit"
$ws.Range("M1").Value = "Comment"

$ws.Range("J2").Value = "# Prise d'essai"
$ws.Range("K2").Value = "# Unité de mesure de la prise d’essai"
$ws.Range("L2").Value = "# Unité du résultat"
$ws.Range("M2").Value = "# Commentaire"

$ws.Range("J3").Value = "#float"
$ws.Range("K3").Value = "#string"
$ws.Range("L3").Value = "#string"
$ws.Range("M3").Value = "#string"

$ws.Range("J4").Value = "# format: nombre décimal, ne pas spécifier d'unité"
$ws.Range("K4").Value = "# format: texte"
$ws.Range("L4").Value = "# format: texte"
$ws.Range("M4").Value = "# format: texte libre"

$ws.Range("J5").Value = "# ex: 2.0"
$ws.Range("K5").Value = "# ex: mg"
$ws.Range("L5").Value = "# ex: mg/ml"
$ws.Range("M5").Value = ""
